# Update "想去人数" (want-to-go count) values in column F across the four
# worksheets (展览, 演出, 本地生活, 全部类型), per the latest scrape refresh.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value  = 193
$ws1.Range("F3").Value  = 406
$ws1.Range("F4").Value  = 1144
$ws1.Range("F7").Value  = 25
$ws1.Range("F8").Value  = 1066
$ws1.Range("F10").Value = 345
$ws1.Range("F11").Value = 421
$ws1.Range("F15").Value = 31
$ws1.Range("F17").Value = 480
$ws1.Range("F19").Value = 5617
$ws1.Range("F21").Value = 1566
$ws1.Range("F22").Value = 369
$ws1.Range("F23").Value = 4798
$ws1.Range("F25").Value = 85
$ws1.Range("F26").Value = 1504
$ws1.Range("F27").Value = 15
$ws1.Range("F30").Value = 69

# --- Sheet 2: 演出 ----------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F5").Value = 138
$ws2.Range("F8").Value = 106

# --- Sheet 3: 本地生活 -------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 9386
$ws3.Range("F4").Value = 2133

# --- Sheet 4: 全部类型 -------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value  = 9386
$ws4.Range("F4").Value  = 2133
$ws4.Range("F5").Value  = 193
$ws4.Range("F6").Value  = 406
$ws4.Range("F7").Value  = 1144
$ws4.Range("F10").Value = 25
$ws4.Range("F11").Value = 1066
$ws4.Range("F12").Value = 345
$ws4.Range("F13").Value = 421
$ws4.Range("F17").Value = 31
$ws4.Range("F23").Value = 5617
$ws4.Range("F25").Value = 1566
$ws4.Range("F28").Value = 369
$ws4.Range("F31").Value = 4798
$ws4.Range("F33").Value = 85
$ws4.Range("F34").Value = 1504
$ws4.Range("F35").Value = 15
$ws4.Range("F38").Value = 69
